$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two "entry" blocks that occupied rows 2-6 are swapped:
#   Block "Publicacion SMQ" (Postdoctoral Research Stays / Orchids Program,
#   Minciencias) used to be first (rows 2-3) and now comes second (rows 5-6).
#   Block "Tutorizacion Postgrado 2023-2024" (Internal Call, Universidad El
#   Bosque) used to be second (rows 4-6) and now comes first (rows 2-4).

# --- New row 2 (previously row 4): Internal Call / Universidad El Bosque ---
$ws.Range("A2").Value = 'IX \href{https://www.unbosque.edu.co/centro-informacion/convocatoria/xiv-convocatoria-interna-de-investigaciones}{Internal Call for Financing Research and Technological Innovation Projects El Bosque University}, 2024'
$ws.Range("B2").Value = 'Jan. 2024 - Jan. 2026'
$ws.Range("C2").Value = '\href{https://www.unbosque.edu.co/}{Universidad El Bosque}'
$ws.Range("D2").Value = 'Bogota, Colombia'
$ws.Range("E2").Value = "Project: Effect of real and simulated resource control on androphilic women's preferences for masculinity in men's faces: an experimental study using eye-tracking"

# --- New row 3 (previously row 5): role line, now plain (non-currency) style ---
$ws.Range("E3").Value = 'Role: Principal Researcher'
$ws.Range("E5").Copy()
$ws.Range("E3").PasteSpecial(-4122)

# --- New row 4 (previously row 6, but A:D now blank): amount line, currency style ---
$ws.Range("A4").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E4").Value = 'COP\$90.000.000'
$ws.Range("E6").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows("4:4").AutoFit()

# --- New row 5 (previously row 2): Postdoctoral Research Stays / Orchids Program ---
$ws.Range("A5").Value = '\href{https://minciencias.gov.co/convocatorias/construccion-paz-programa-y-proyectos-ctei-fortalecimiento-capacidades-para-la}{Postdoctoral Research Stays -  Call 935-2023 - Orchids Program. Women in science: agents for peace: Agents for Peace 2023}'
$ws.Range("B5").Value = 'Dic. 2023 - Jan. 2025'
$ws.Range("C5").Value = '\href{https://minciencias.gov.co/}{Minciencias}'
$ws.Range("D5").Value = 'Barranquilla, Colombia'
$ws.Range("E5").Value = "Project: Effect of resource availability on women's preferences for masculinity faces in interaction with hormonal, cognitive, and socio-contextual factors such as actual resource scarcity and exposure to violence: an experimental study using eye-tracking"
$ws.Rows("5:5").RowHeight = 43.2

# --- New row 6 (previously row 3): amount line, stays currency style ---
$ws.Range("E6").Value = 'COP\$356.040.884 '

$excel.CutCopyMode = 0

# Selection moves to A7:XFD9 (active cell A7), matching the post-edit cursor position
$ws.Range("A7:XFD9").Select()
